$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.300369739532471
$ws.Range("B1").Value = 4.731704235076904
$ws.Range("C1").Value = 7.632800102233887
$ws.Range("D1").Value = 7.336466312408447
$ws.Range("E1").Value = 4.782388687133789
